# Updated symbol list on Sun Jan 15 23:56:23 UTC 2023 with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for the
# crypto table on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Address,
        [string]$NewValue
    )
    $cell = $ws.Range($Address)
    # Force text number-format so numeric-/percent-looking strings are
    # stored verbatim instead of being auto-converted by Excel into a
    # real number.
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    # Put formatting back to the default so we do not leave a stray
    # number-format behind on cells that originally had none.
    $cell.Style = "Normal"
}

Set-TextCell "D2" "301.55"
Set-TextCell "E2" "-1.14%"
Set-TextCell "D3" "31.29"
Set-TextCell "E3" "-3.33%"
Set-TextCell "D4" "5.129"
Set-TextCell "E4" "-3.54%"
Set-TextCell "D5" "0.07405"
Set-TextCell "E5" "-2.45%"
Set-TextCell "D6" "2.191"
Set-TextCell "E6" "14.85%"
Set-TextCell "D7" "7.934"
Set-TextCell "E7" "0.43%"
Set-TextCell "D8" "3.827"
Set-TextCell "E8" "-1.27%"
Set-TextCell "D9" "0.9200"
Set-TextCell "E9" "-1.32%"
Set-TextCell "E10" "-0.25%"
Set-TextCell "D11" "0.07585"
Set-TextCell "E11" "-4.95%"
Set-TextCell "D12" "0.08114"
Set-TextCell "E12" "0.11%"
Set-TextCell "D13" "0.02997"
Set-TextCell "E13" "-2.02%"
Set-TextCell "D14" "0.09917"
Set-TextCell "E14" "-0.32%"
Set-TextCell "D15" "0.001496"
Set-TextCell "E15" "-0.31%"
Set-TextCell "D16" "0.006111"
Set-TextCell "E16" "-1.67%"
Set-TextCell "D17" "3.476"
Set-TextCell "E17" "0.90%"
Set-TextCell "E18" "-0.25%"
Set-TextCell "D19" "0.3262"
Set-TextCell "E19" "-1.18%"
Set-TextCell "D20" "0.1319"
Set-TextCell "E20" "-1.89%"
Set-TextCell "D21" "4.644"
Set-TextCell "E21" "1.85%"
Set-TextCell "D22" "0.04654"
Set-TextCell "E22" "1.11%"
Set-TextCell "E23" "-3.18%"
Set-TextCell "D24" "0.001225"
Set-TextCell "E24" "0.90%"
Set-TextCell "D25" "0.004481"
Set-TextCell "E26" "-7.07%"
Set-TextCell "D27" "0.0003428"
Set-TextCell "E27" "92.35%"
Set-TextCell "D39" "0.01722"
Set-TextCell "E39" "-2.00%"
Set-TextCell "D40" "0.04505"
Set-TextCell "E40" "-1.47%"
Set-TextCell "D41" "0.007333"
Set-TextCell "E41" "4.13%"
Set-TextCell "D42" "0.1347"
Set-TextCell "E42" "-1.28%"
Set-TextCell "E43" "7.81%"
Set-TextCell "D44" "0.01061"
Set-TextCell "E44" "-23.17%"
Set-TextCell "D45" "0.00006280"
Set-TextCell "E45" "2.05%"
Set-TextCell "E46" "13.42%"
Set-TextCell "E47" "-18.21%"
